$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove empty placeholder "imagem" cells (column R) for rows 2-111 ---
# Every existing data row had an empty inline-string cell in column R ("imagem").
# The edit drops those empty placeholders for rows 2-111 (column R = 18th column).
for ($r = 2; $r -le 111; $r++) {
    $ws.Cells.Item($r, 18).ClearContents()
}

# --- Append new question rows 113-123 ---

# Row 113
$ws.Range('A113').Value = 113
$ws.Range('B113').Value = 'CESGRANRIO'
$ws.Range('C113').Value = 'UNIRIO'
$ws.Range('D113').NumberFormat = "@"
$ws.Range('D113').Value = '2016'
$ws.Range('E113').Value = 'Em “Fazia calor no Rio, 40 graus e qualquer coisa, quase 41.” (l. 1-2), o uso do pretérito imperfeito do indicativo busca.'
$ws.Range('F113').Value = 'Português'
$ws.Range('G113').Value = 'Emprego de Tempos e Modos'
$ws.Range('H113').Value = 'Médio'
$ws.Range('I113').Value = 'ME'
$ws.Range('J113').Value = 'estabelecer uma relação de causa e efeito'
$ws.Range('K113').Value = 'contextualizar o tempo da narrativa'
$ws.Range('L113').Value = 'introduzir uma ambiência de suspense'
$ws.Range('M113').Value = 'banalizar o calor que fazia no Rio'
$ws.Range('N113').Value = 'projetar uma possibilidade'
$ws.Range('O113').Value = 'B'
$ws.Range('P113').Value = 0
$ws.Range('Q113').Value = 0
$ws.Range('R113').Value = '7b209d1e-87f0-48a2-b97f-00b3f541acb2.png'

# Row 114
$ws.Range('A114').Value = 114
$ws.Range('B114').Value = 'CESGRANRIO'
$ws.Range('C114').Value = 'BANCO DA AMAZÔNIA'
$ws.Range('D114').NumberFormat = "@"
$ws.Range('D114').Value = '2015'
$ws.Range('E114').Value = 'O verbo em destaque está conjugado de acordo com a norma-padrão em:'
$ws.Range('F114').Value = 'Português'
$ws.Range('G114').Value = 'Emprego de Tempos e Modos'
$ws.Range('H114').Value = 'Médio'
$ws.Range('I114').Value = 'ME'
$ws.Range('J114').Value = '<u>Pegue </u>o outro elevador, por favor'
$ws.Range('K114').Value = 'É preciso que você <u>esteje </u>atento a situações de perigo'
$ws.Range('L114').Value = 'Será muito bom se você <u>propor </u>um outro acesso aos passageiros'
$ws.Range('M114').Value = '<u>Seje</u> sempre bem-humorado com os passageiros'
$ws.Range('N114').Value = 'Gostaríamos de que você <u>vesse </u>esse filme'
$ws.Range('O114').Value = 'A'
$ws.Range('P114').Value = 0
$ws.Range('Q114').Value = 0

# Row 115
$ws.Range('A115').Value = 115
$ws.Range('B115').Value = 'CESGRANRIO'
$ws.Range('C115').Value = 'LIQUIGÁS'
$ws.Range('D115').NumberFormat = "@"
$ws.Range('D115').Value = '2015'
$ws.Range('E115').Value = 'A forma verbal destacada está empregada de acordo com a norma-padrão em:'
$ws.Range('F115').Value = 'Português'
$ws.Range('G115').Value = 'Emprego de Tempos e Modos'
$ws.Range('H115').Value = 'Médio'
$ws.Range('I115').Value = 'ME'
$ws.Range('J115').Value = 'Quando as pessoas <u>fazerem </u>compras nas lojas locais, poderão usar o cartão de crédito comunitário'
$ws.Range('K115').Value = 'Os consumidores preocupados com os gastos tinham <u>trago </u>pouco dinheiro para as suas compras'
$ws.Range('L115').Value = 'Os financiamentos serão ampliados quando os bancos <u>estarem </u>com os juros baixos'
$ws.Range('M115').Value = 'O ideal seria que os clientes dos bancos comunitários <u>pudessem </u>aumentar sua renda mensal'
$ws.Range('N115').Value = 'Se os bancos <u>darem </u>mais crédito aos moradores, aumentará a construção de casas na comunidade'
$ws.Range('O115').Value = 'D'
$ws.Range('P115').Value = 0
$ws.Range('Q115').Value = 0

# Row 116
$ws.Range('A116').Value = 116
$ws.Range('B116').Value = 'CESGRANRIO'
$ws.Range('C116').Value = 'LIQUIGÁS'
$ws.Range('D116').NumberFormat = "@"
$ws.Range('D116').Value = '2015'
$ws.Range('E116').Value = 'Em “os movimentos caóticos atuais já <u>eram </u>os primeiros passos afinando-se e orquestrando-se para uma situação econômica mais digna", observa-se a adequada flexão do verbo destacado, o que também se verifica em'
$ws.Range('F116').Value = 'Português'
$ws.Range('G116').Value = 'Emprego de Tempos e Modos'
$ws.Range('H116').Value = 'Médio'
$ws.Range('I116').Value = 'ME'
$ws.Range('J116').Value = 'O povo brasileiro sois articulado politicamente'
$ws.Range('K116').Value = 'Se fôssemos menos explorados, o país prosperaria mais'
$ws.Range('L116').Value = 'Quando fores corretos os políticos brasileiros, a fome terá fim'
$ws.Range('M116').Value = 'É fundamental que sejemos um povo mais politizado'
$ws.Range('N116').Value = 'A população deseja que o sistema de saúde sejais melhor'
$ws.Range('O116').Value = 'B'
$ws.Range('P116').Value = 0
$ws.Range('Q116').Value = 0

# Row 117
$ws.Range('A117').Value = 117
$ws.Range('B117').Value = 'CESGRANRIO'
$ws.Range('C117').Value = 'BANCO DO BRASIL'
$ws.Range('D117').NumberFormat = "@"
$ws.Range('D117').Value = '2015'
$ws.Range('E117').Value = 'O emprego do verbo destacado no trecho “queremos contribuir para o crescimento sustentável das empresas” contribui para indicar uma pretensão do presidente do Sindicato dos Lojistas, que começa no presente e se estende no futuro.
Se, respeitando-se o contexto original, a frase indicasse uma pretensão que começasse no passado e se estendesse no tempo, o verbo adequado seria o que se destaca em:'
$ws.Range('F117').Value = 'Português'
$ws.Range('G117').Value = 'Emprego de Tempos e Modos'
$ws.Range('H117').Value = 'Médio'
$ws.Range('I117').Value = 'ME'
$ws.Range('J117').Value = '<u>quisemos </u>contribuir para o crescimento sustentável das empresas'
$ws.Range('K117').Value = '<u>quisermos </u>contribuir para o crescimento sustentável das empresas'
$ws.Range('L117').Value = '<u>quiséssemos </u>contribuir para o crescimento sustentável das empresas'
$ws.Range('M117').Value = '<u>quereremos </u>contribuir para o crescimento sustentável das empresas'
$ws.Range('N117').Value = '<u>quisera </u>poder contribuir para o crescimento sustentável das empresas'
$ws.Range('O117').Value = 'A'
$ws.Range('P117').Value = 0
$ws.Range('Q117').Value = 0

# Row 118
$ws.Range('A118').Value = 118
$ws.Range('B118').Value = 'CESGRANRIO'
$ws.Range('C118').Value = 'LIQUIGÁS'
$ws.Range('D118').NumberFormat = "@"
$ws.Range('D118').Value = '2014'
$ws.Range('E118').Value = 'Flexionado na 1ª pessoa do singular do presente do indicativo, o verbo fazer assume forma irregular: faço.
O mesmo acontece com o seguinte verbo:'
$ws.Range('F118').Value = 'Português'
$ws.Range('G118').Value = 'Emprego de Tempos e Modos'
$ws.Range('H118').Value = 'Médio'
$ws.Range('I118').Value = 'ME'
$ws.Range('J118').Value = 'depender'
$ws.Range('K118').Value = 'dominar'
$ws.Range('L118').Value = 'medir'
$ws.Range('M118').Value = 'pensar'
$ws.Range('N118').Value = 'dever'
$ws.Range('O118').Value = 'C'
$ws.Range('P118').Value = 0
$ws.Range('Q118').Value = 0

# Row 119
$ws.Range('A119').Value = 119
$ws.Range('B119').Value = 'CESGRANRIO'
$ws.Range('C119').Value = 'LIQUIGÁS'
$ws.Range('D119').NumberFormat = "@"
$ws.Range('D119').Value = '2014'
$ws.Range('E119').Value = 'A frase em que a flexão do verbo auxiliar destacado obedece aos princípios da norma padrão é'
$ws.Range('F119').Value = 'Português'
$ws.Range('G119').Value = 'Emprego de Tempos e Modos'
$ws.Range('H119').Value = 'Médio'
$ws.Range('I119').Value = 'ME'
$ws.Range('J119').Value = 'Alguns estudiosos consideram que <u>podem </u>haver robôs tão inteligentes quanto o homem'
$ws.Range('K119').Value = '<u>Devem </u>existir formas de garantir a exploração de outras tarefas destinadas aos robôs'
$ws.Range('L119').Value = 'No futuro, <u>devem </u>haver outras formas de investimentos para garantir a evolução da robótica'
$ws.Range('M119').Value = '<u>Pode </u>existir obstáculos que os robôs sejam capazes de superar, como a locomoção e o diálogo'
$ws.Range('N119').Value = '<u>Pode </u>surgir novas tecnologias para aperfeiçoar a conquista espacial'
$ws.Range('O119').Value = 'B'
$ws.Range('P119').Value = 0
$ws.Range('Q119').Value = 0

# Row 120
$ws.Range('A120').Value = 120
$ws.Range('B120').Value = 'CESGRANRIO'
$ws.Range('C120').Value = 'LIQUIGÁS'
$ws.Range('D120').NumberFormat = "@"
$ws.Range('D120').Value = '2014'
$ws.Range('E120').Value = 'No trecho “Para o futuro, prenunciam-se perguntas mais difíceis, mais desafiadoras — e até ameaçadoras — do que aquelas relativas ao uso de drones.”, o verbo prenunciar foi utilizado no plural por se tratar de uma construção de passiva pronominal com sujeito também no plural. O mesmo procedimento é adotado no verbo destacado em:'
$ws.Range('F120').Value = 'Português'
$ws.Range('G120').Value = 'Emprego de Tempos e Modos'
$ws.Range('H120').Value = 'Médio'
$ws.Range('I120').Value = 'ME'
$ws.Range('J120').Value = 'Para conquistar posição de vanguarda na corrida espacial, <u>obedecem</u>-se a princípios básicos de inovação tecnológica'
$ws.Range('K120').Value = 'Na missão espacial ao solo marciano, ambiente inóspito aos humanos, <u>assistiram</u>-se a episódios inesquecíveis'
$ws.Range('L120').Value = 'Nos livros e filmes de ficção científica do século passado, <u>falavam</u>-se de robôs como uma possibilidade muito próxima e viável'
$ws.Range('M120').Value = 'Com o avanço das pesquisas em robótica nas últimas décadas, <u>delegam</u>-se atividades eminentemente humanas às máquinas'
$ws.Range('N120').Value = 'Para evitar que o crescimento da robótica provoque distorções incontroláveis, <u>necessitam</u>-se de leis protecionistas'
$ws.Range('O120').Value = 'D'
$ws.Range('P120').Value = 0
$ws.Range('Q120').Value = 0

# Row 121
$ws.Range('A121').Value = 121
$ws.Range('B121').Value = 'CESGRANRIO'
$ws.Range('C121').Value = 'ELETRONUCLEAR'
$ws.Range('D121').NumberFormat = "@"
$ws.Range('D121').Value = '2022'
$ws.Range('E121').Value = 'Em que frase o verbo irregular destacado está empregado de acordo com a norma-padrão da Língua Portuguesa?'
$ws.Range('F121').Value = 'Português'
$ws.Range('G121').Value = 'Verbos Traiçoeiros'
$ws.Range('H121').Value = 'Médio'
$ws.Range('I121').Value = 'ME'
$ws.Range('J121').Value = 'Os médicos <u>preveram </u>que ela teria complicações da doença. (verbo PREVER)'
$ws.Range('K121').Value = 'Se eu me <u>oposse </u>a suas orientações, ela me advertia. (verbo OPOR)'
$ws.Range('L121').Value = 'Minha mãe sempre me <u>acodia </u>nos momentos difíceis. (verbo ACUDIR)'
$ws.Range('M121').Value = 'Maria José sempre <u>soube </u>defender filhos e netos. (verbo SABER)'
$ws.Range('N121').Value = 'Quando entrava numa briga, ela sempre <u>intervia </u>em meu favor. (verbo INTERVIR)'
$ws.Range('O121').Value = 'D'
$ws.Range('P121').Value = 0
$ws.Range('Q121').Value = 0

# Row 122
$ws.Range('A122').Value = 122
$ws.Range('B122').Value = 'CESGRANRIO'
$ws.Range('C122').Value = 'LIQUIGÁS'
$ws.Range('D122').NumberFormat = "@"
$ws.Range('D122').Value = '2014'
$ws.Range('E122').Value = 'A forma verbal destacada está empregada de acordo com a norma-padrão da Língua Portuguesa em:'
$ws.Range('F122').Value = 'Português'
$ws.Range('G122').Value = 'Verbos Traiçoeiros'
$ws.Range('H122').Value = 'Médio'
$ws.Range('I122').Value = 'ME'
$ws.Range('J122').Value = 'Se os governantes <u>verem </u>o prejuízo causado pelas variações do clima, talvez tomem medidas cautelares'
$ws.Range('K122').Value = 'A construção de novas hidrelétricas dependia de que as verbas se <u>mantessem </u>inalteradas'
$ws.Range('L122').Value = 'As variações do clima não afetariam o meio ambiente se a população <u>interviesse </u>nas políticas públicas'
$ws.Range('M122').Value = 'Todos <u>ansiam </u>que os eventos climáticos extremos não cheguem a causar problemas ambientais'
$ws.Range('N122').Value = 'Um grupo de pesquisadores <u>entreveu </u>a possibilidade de prejuízos na produção de energia por causa das alterações das chuvas na Amazônia'
$ws.Range('O122').Value = 'C'
$ws.Range('P122').Value = 0
$ws.Range('Q122').Value = 0

# Row 123
$ws.Range('A123').Value = 123
$ws.Range('B123').Value = 'CESGRANRIO'
$ws.Range('C123').Value = 'LIQUIGÁS'
$ws.Range('D123').NumberFormat = "@"
$ws.Range('D123').Value = '2014'
$ws.Range('E123').Value = 'A forma verbal destacada está empregada de acordo com a norma-padrão da Língua Portuguesa em:'
$ws.Range('F123').Value = 'Português'
$ws.Range('G123').Value = 'Verbos Traiçoeiros'
$ws.Range('H123').Value = 'Médio'
$ws.Range('I123').Value = 'ME'
$ws.Range('J123').Value = 'Se os governantes <u>verem </u>o prejuízo causado pelas variações do clima, talvez tomem medidas cautelares'
$ws.Range('K123').Value = 'A construção de novas hidrelétricas dependia de que as verbas se <u>mantessem </u>inalteradas'
$ws.Range('L123').Value = 'As variações do clima não afetariam o meio ambiente se a população <u>interviesse </u>nas políticas públicas'
$ws.Range('M123').Value = 'Todos <u>ansiam </u>que os eventos climáticos extremos não cheguem a causar problemas ambientais'
$ws.Range('N123').Value = 'Um grupo de pesquisadores <u>entreveu </u>a possibilidade de prejuízos na produção de energia por causa das alterações das chuvas na Amazônia.G ABARITO 1. LETRA D 2. LETRA C'
$ws.Range('O123').Value = 'C'
$ws.Range('P123').Value = 0
$ws.Range('Q123').Value = 0
# R123: empty placeholder (imagem) cell - intentionally left blank
